$wb = $excel.ActiveWorkbook

# --- Sheet 1: Neg_Change ---
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws1.Range("A2").Value = "SBILIFE"
$ws1.Range("B2").Value = 1967
$ws1.Range("C2").Value = 1977
$ws1.Range("D2").Value = 1951.2
$ws1.Range("E2").Value = 1974
$ws1.Range("F2").Value = 483716
$ws1.Range("G2").Value = 1011160
$ws1.Range("H2").Value = -0.5216226907709957
$ws1.Range("I2").Value = "SBILIFE"

$ws1.Range("A3").Value = "TATACONSUM"
$ws1.Range("B3").Value = 1172.4
$ws1.Range("C3").Value = 1178
$ws1.Range("D3").Value = 1159.6
$ws1.Range("E3").Value = 1162.5
$ws1.Range("F3").Value = 487062
$ws1.Range("G3").Value = 1020242
$ws1.Range("H3").Value = -0.5226015004283298
$ws1.Range("I3").Value = "TATACONSUM"

$ws1.Range("A4").Value = "ADANIENT"
$ws1.Range("B4").Value = 2300
$ws1.Range("C4").Value = 2304
$ws1.Range("D4").Value = 2257
$ws1.Range("E4").Value = 2259.9
$ws1.Range("F4").Value = 1313265
$ws1.Range("G4").Value = 3037620
$ws1.Range("H4").Value = -0.5676664625595038
$ws1.Range("I4").Value = "ADANIENT"

$ws1.Range("A5").Value = "DIVISLAB"
$ws1.Range("B5").Value = 6500.5
$ws1.Range("C5").Value = 6524.5
$ws1.Range("D5").Value = 6404.5
$ws1.Range("E5").Value = 6419.5
$ws1.Range("F5").Value = 100192
$ws1.Range("G5").Value = 250130
$ws1.Range("H5").Value = -0.5994402910486547
$ws1.Range("I5").Value = "DIVISLAB"

$ws1.Range("A6").Value = "GODREJCP"
$ws1.Range("B6").Value = 1145
$ws1.Range("C6").Value = 1145.9
$ws1.Range("D6").Value = 1122.7
$ws1.Range("E6").Value = 1129
$ws1.Range("F6").Value = 766007
$ws1.Range("G6").Value = 1513316
$ws1.Range("H6").Value = -0.4938221759368169
$ws1.Range("I6").Value = "GODREJCP"

$ws1.Range("A7").Value = "LODHA"
$ws1.Range("B7").Value = 1156.3
$ws1.Range("C7").Value = 1158
$ws1.Range("D7").Value = 1125.1
$ws1.Range("E7").Value = 1126.6
$ws1.Range("F7").Value = 684453
$ws1.Range("G7").Value = 1536998
$ws1.Range("H7").Value = -0.5546819189094586
$ws1.Range("I7").Value = "LODHA"

$ws1.Range("A8").Value = "360ONE"
$ws1.Range("B8").Value = 1184.9
$ws1.Range("C8").Value = 1198.6
$ws1.Range("D8").Value = 1181.2
$ws1.Range("E8").Value = 1190
$ws1.Range("F8").Value = 560894
$ws1.Range("G8").Value = 1176162
$ws1.Range("H8").Value = -0.5231150130679277
$ws1.Range("I8").Value = "360ONE"

$ws1.Range("A9").Value = "LICHSGFIN"
$ws1.Range("B9").Value = 550.9
$ws1.Range("C9").Value = 554.75
$ws1.Range("D9").Value = 549.05
$ws1.Range("E9").Value = 550
$ws1.Range("F9").Value = 289709
$ws1.Range("G9").Value = 618499
$ws1.Range("H9").Value = -0.5315934221397286
$ws1.Range("I9").Value = "LICHSGFIN"

$ws1.Range("A10").Value = "BHARATFORG"
$ws1.Range("B10").Value = 1438
$ws1.Range("C10").Value = 1450
$ws1.Range("D10").Value = 1424.9
$ws1.Range("E10").Value = 1430.4
$ws1.Range("F10").Value = 495952
$ws1.Range("G10").Value = 1148404
$ws1.Range("H10").Value = -0.5681380420130895
$ws1.Range("I10").Value = "BHARATFORG"

$ws1.Range("A11").Value = "MCX"
$ws1.Range("B11").Value = 10200
$ws1.Range("C11").Value = 10308
$ws1.Range("D11").Value = 10137
$ws1.Range("E11").Value = 10214
$ws1.Range("F11").Value = 427516
$ws1.Range("G11").Value = 918982
$ws1.Range("H11").Value = -0.5347939350281072
$ws1.Range("I11").Value = "MCX"

$ws1.Range("A12").Value = "PGEL"
$ws1.Range("B12").Value = 590
$ws1.Range("C12").Value = 597.65
$ws1.Range("D12").Value = 586.55
$ws1.Range("E12").Value = 592.55
$ws1.Range("F12").Value = 930540
$ws1.Range("G12").Value = 2063578
$ws1.Range("H12").Value = -0.5490647797175585
$ws1.Range("I12").Value = "PGEL"

# --- Sheet 2: Pos_Change ---
$ws2 = $wb.Worksheets.Item("Pos_Change")
$ws2.Range("A2").Value = "ULTRACEMCO"
$ws2.Range("B2").Value = 11620
$ws2.Range("C2").Value = 12021
$ws2.Range("D2").Value = 11556
$ws2.Range("E2").Value = 12013
$ws2.Range("F2").Value = 228729
$ws2.Range("G2").Value = 159990
$ws2.Range("H2").Value = 0.4296456028501781
$ws2.Range("I2").Value = "ULTRACEMCO"

$ws2.Range("A3").Value = "BEL"
$ws2.Range("B3").Value = 415.9
$ws2.Range("C3").Value = 419.2
$ws2.Range("D3").Value = 413.5
$ws2.Range("E3").Value = 417
$ws2.Range("F3").Value = 9650320
$ws2.Range("G3").Value = 6690232
$ws2.Range("H3").Value = 0.4424492304601694
$ws2.Range("I3").Value = "BEL"

$ws2.Range("A4").Value = "ADANIPORTS"
$ws2.Range("B4").Value = 1530
$ws2.Range("C4").Value = 1549
$ws2.Range("D4").Value = 1523
$ws2.Range("E4").Value = 1531.5
$ws2.Range("F4").Value = 4059126
$ws2.Range("G4").Value = 2612938
$ws2.Range("H4").Value = 0.5534719920641056
$ws2.Range("I4").Value = "ADANIPORTS"

$ws2.Range("A5").Value = "NAUKRI"
$ws2.Range("B5").Value = 1326.4
$ws2.Range("C5").Value = 1369
$ws2.Range("D5").Value = 1326.4
$ws2.Range("E5").Value = 1361
$ws2.Range("F5").Value = 939714
$ws2.Range("G5").Value = 655737
$ws2.Range("H5").Value = 0.4330653905452949
$ws2.Range("I5").Value = "NAUKRI"

$ws2.Range("A6").Value = "JSWENERGY"
$ws2.Range("B6").Value = 488.7
$ws2.Range("C6").Value = 493.4
$ws2.Range("D6").Value = 484.1
$ws2.Range("E6").Value = 491.45
$ws2.Range("F6").Value = 1933563
$ws2.Range("G6").Value = 1256225
$ws2.Range("H6").Value = 0.5391852574180581
$ws2.Range("I6").Value = "JSWENERGY"

$ws2.Range("A7").Value = "NMDC"
$ws2.Range("B7").Value = 74.2
$ws2.Range("C7").Value = 75.75
$ws2.Range("D7").Value = 73.87
$ws2.Range("E7").Value = 75.45
$ws2.Range("F7").Value = 24997800
$ws2.Range("G7").Value = 16024318
$ws2.Range("H7").Value = 0.55999150790692
$ws2.Range("I7").Value = "NMDC"

$ws2.Range("A8").Value = "IDFCFIRSTB"
$ws2.Range("B8").Value = 80.48999999999999
$ws2.Range("C8").Value = 81.04000000000001
$ws2.Range("D8").Value = 79.95999999999999
$ws2.Range("E8").Value = 80.55
$ws2.Range("F8").Value = 21543164
$ws2.Range("G8").Value = 14855394
$ws2.Range("H8").Value = 0.4501913581019796
$ws2.Range("I8").Value = "IDFCFIRSTB"

$ws2.Range("A9").Value = "BANKINDIA"
$ws2.Range("B9").Value = 147.5
$ws2.Range("C9").Value = 149.04
$ws2.Range("D9").Value = 146.35
$ws2.Range("E9").Value = 147.35
$ws2.Range("F9").Value = 6729070
$ws2.Range("G9").Value = 4416863
$ws2.Range("H9").Value = 0.5234952951902742
$ws2.Range("I9").Value = "BANKINDIA"

$ws2.Range("A10").Value = "UNOMINDA"
$ws2.Range("B10").Value = 1307
$ws2.Range("C10").Value = 1337.5
$ws2.Range("D10").Value = 1302
$ws2.Range("E10").Value = 1307
$ws2.Range("F10").Value = 884017
$ws2.Range("G10").Value = 582383
$ws2.Range("H10").Value = 0.517930640145746
$ws2.Range("I10").Value = "UNOMINDA"

$ws2.Range("A11").Value = "PATANJALI"
$ws2.Range("B11").Value = 568.15
$ws2.Range("C11").Value = 570.8
$ws2.Range("D11").Value = 562.8
$ws2.Range("E11").Value = 567.85
$ws2.Range("F11").Value = 553676
$ws2.Range("G11").Value = 360701
$ws2.Range("H11").Value = 0.5349999029667232
$ws2.Range("I11").Value = "PATANJALI"

$ws2.Range("A12").Value = "GMRAIRPORT"
$ws2.Range("B12").Value = 108.9
$ws2.Range("C12").Value = 109.6
$ws2.Range("D12").Value = 106.83
$ws2.Range("E12").Value = 107.67
$ws2.Range("F12").Value = 40681263
$ws2.Range("G12").Value = 25821928
$ws2.Range("H12").Value = 0.5754541256563026
$ws2.Range("I12").Value = "GMRAIRPORT"

$ws2.Range("A13").Value = "GRANULES"
$ws2.Range("B13").Value = 562.5
$ws2.Range("C13").Value = 571.1
$ws2.Range("D13").Value = 559.05
$ws2.Range("E13").Value = 568.5
$ws2.Range("F13").Value = 1311056
$ws2.Range("G13").Value = 857235
$ws2.Range("H13").Value = 0.5294009227341394
$ws2.Range("I13").Value = "GRANULES"

$ws2.Range("A14").Value = "BANDHANBNK"
$ws2.Range("B14").Value = 151.11
$ws2.Range("C14").Value = 151.2
$ws2.Range("D14").Value = 149.38
$ws2.Range("E14").Value = 150.2
$ws2.Range("F14").Value = 6906158
$ws2.Range("G14").Value = 4578914
$ws2.Range("H14").Value = 0.5082523934714651
$ws2.Range("I14").Value = "BANDHANBNK"

